$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("股票")

# Insert a new "category" column between "property_category" (H) and "date" (old I).
# This shifts date/legislator_name/legislator_id from I/J/K to J/K/L.
$ws.Columns.Item(9).Insert()

# Append two new trailing columns: "source_file" and "index".
$ws.Columns.Item(13).Insert()
$ws.Columns.Item(14).Insert()

$ws.Cells.Item(1, 9).Value = "category"
$ws.Cells.Item(1, 13).Value = "source_file"
$ws.Cells.Item(1, 14).Value = "index"

$lastRow = 13
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 9).Value = "normal"
    $ws.Cells.Item($r, 13).Value = "tmp4cfc1"
    $ws.Cells.Item($r, 14).Value = $ws.Cells.Item($r, 1).Value()
}
